# REPORTGEN-1070: update page title in templates
#
# 1) Bump the cached "datetimeFigureOut" field text from 10/22/2021 to
#    10/25/2021 everywhere it appears (slide master, the three slide
#    layouts that carry their own date placeholder, and the notes master).
# 2) Drop the " - sample 2" suffix from two chart-demo slide titles.

$p = $ppt.ActivePresentation

$oldDate = "10/22/2021"
$newDate = "10/25/2021"

# --- Slide master: "Date Placeholder 3" ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide layouts that define their own date placeholder ---
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Notes master: "Date Placeholder 2" ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide titles: drop the " <en-dash> sample 2" suffix ---
$dash = [char]0x2013
$clusteredOld = "Clustered column graph $dash sample 2"
$stackedOld = "Stacked Bar $dash sample 2"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = $slide.Shapes.Item(1)
    if ($title.HasTextFrame) {
        $text = $title.TextFrame.TextRange.Text
        if ($text -eq $clusteredOld) {
            $title.TextFrame.TextRange.Text = "Clustered column graph"
        } elseif ($text -eq $stackedOld) {
            $title.TextFrame.TextRange.Text = "Stacked Bar"
        }
    }
}
